$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BalancoResumido")

# Data rows for " 01/02/03 Junho de 2016" as already present earlier in the
# sheet (e.g. rows 180/182/183) - reused here to extend the table so the
# balanco detalhado block can be printed.
$jun01 = @(" 01 Junho de 2016","39.365","11.382","1.990","9.025","2.761","64.523","1.042","63.481","38.945","11.235","2.016","8.363","3.377","63.936","1.003","62.933")
$jun02 = @(" 02 Junho de 2016","38.527","11.318","1.990","9.623","3.070","64.528","1.074","63.454","39.694","11.242","2.014","8.743","2.861","64.554","925","63.629")
$jun03 = @(" 03 Junho de 2016","37.840","10.776","1.990","9.849","2.564","63.019","0","63.019","38.549","11.075","2.014","9.261","2.204","63.103","0","63.103")

# Rows 303..310 repeat that three-day block (with the last partial repeat
# trailing off after a single day), matching how it was pasted in originally.
$rows = @($jun01, $jun02, $jun03, $jun01, $jun02, $jun01, $jun02, $jun01)

$startRow = 303
for ($r = 0; $r -lt $rows.Length; $r++) {
    $targetRow = $startRow + $r
    $vals = $rows[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $cell = $ws.Cells.Item($targetRow, $c + 1)
        # Force text type so numeric-looking strings ("39.365") aren't
        # reinterpreted as numbers, matching the original text cells.
        $cell.NumberFormat = "@"
        $cell.Value = $vals[$c]
        # Columns with a column-level default style (A, D, E, H, L, M) would
        # otherwise pick up an explicit style on write; reset to Normal so
        # the cell matches the plain/default-styled cells used elsewhere in
        # this table.
        $cell.Style = "Normal"
    }
}

# The stray blank "R" cell that used to trail the last row of the table now
# trails the new last row instead (it carries no real value, Excel stores it
# as an empty text cell, i.e. shared-string index 0 = "Programado").
$ws.Range("R302").ClearContents()
$rLast = $ws.Cells.Item(310, 18)
$rLast.NumberFormat = "@"
$rLast.Value = "Programado"
$rLast.Style = "Normal"
